# Weekly fruit/vegetable price update: a new weekly record is inserted
# at row 137 (pushing the existing rows 137-177 down to 138-178), extending
# the table from A1:R177 to A1:R178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 137; everything below (old rows 137-177) shifts
# down to 138-178, and row formatting (e.g. the date style on column D) is
# carried along automatically by Excel's Insert.
$ws.Rows("137:137").Insert()

# Populate the newly inserted row with this week's record (same
# market/category boilerplate as every other row in this block).
$ws.Range("A137").Value = 3
$ws.Range("B137").Value = "Femacal de La Calera"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44463
$ws.Range("E137").Value = 5
$ws.Range("F137").Value = 100112039
$ws.Range("G137").Value = "Ciboulette"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 160
$ws.Range("K137").Value = 1500
$ws.Range("L137").Value = 1500
$ws.Range("M137").Value = 1500
$ws.Range("N137").Value = "`$/docena de atados"
$ws.Range("O137").Value = "Provincia de Quillota"
$ws.Range("P137").Value = 500
$ws.Range("Q137").Value = 3
$ws.Range("R137").Value = "Hortaliza"
